# Generate Report for Handback
# Update handback timestamps / priority values produced by a new report run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# "Correspond Handback DateTime" (col G) for the 2ba3279f... and 86771222... rows
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-03 14:19:24"
$wsOverview.Range("G4").Value = "2016-09-03 14:19:24"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority column (E): ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
# Correspond Handoff Datetime (H)
$wsZhCn.Range("H2").Value = "2016-09-03 14:19:19"
$wsZhCn.Range("H4").Value = "2016-09-03 14:19:19"
# Correspond Handback DateTime (K)
$wsZhCn.Range("K2").Value = "2016-09-03 14:19:37"
$wsZhCn.Range("K4").Value = "2016-09-03 14:19:37"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Priority column (E): ht -> mt
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
# Correspond Handback DateTime (K)
$wsDeDe.Range("K2").Value = "2016-09-03 14:19:43"
$wsDeDe.Range("K4").Value = "2016-09-03 14:19:43"

Write-Output "done"
